# Append: 2025-12-11 06:30 JST
# The scraper re-ran and produced a fresh top-N list (N shrank from 15 to 8
# data rows). Existing rows 2-9 are overwritten in place with the new
# data, and old rows 10-16 are removed entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Drop the now-gone tail rows (old rows 10-16) BEFORE touching rows
#    2-9 so row indices used below stay valid.
# ---------------------------------------------------------------------
$ws.Rows("10:16").Delete()

# ---------------------------------------------------------------------
# 2. Column width tweaks (B, D, H got narrower).
#    Excel's `ColumnWidth` is expressed in "characters of the workbook's
#    default font" and is ~0.8333333333333 wider than the raw OOXML
#    `width` attribute for this font/theme, so we subtract that constant
#    offset to land exactly on the target raw width.
# ---------------------------------------------------------------------
$charOffset = 0.8333333333333
$ws.Columns("B").ColumnWidth = 38 - $charOffset
$ws.Columns("D").ColumnWidth = 28 - $charOffset
$ws.Columns("H").ColumnWidth = 12 - $charOffset

# ---------------------------------------------------------------------
# 3. Refresh the data rows (2-9) with the new scrape.
# ---------------------------------------------------------------------

# Row 2
$ws.Range("A2").Value = "2025-12-11 06:30:42"
$ws.Range("B2").Value = "【急募】AIチャットボット開発のプロフェッショナルを探しています!"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5451734"
$ws.Range("G2").Value = 368
$ws.Range("H2").Value = "🔥AI,Ai ◆開発"

# Row 3
$ws.Range("A3").Value = "2025-12-11 06:30:42"
$ws.Range("B3").Value = "【募集】RPAツール「RoboTANGO」設定代行の専門家を探しています"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5405023"
$ws.Range("G3").Value = 178
$ws.Range("H3").Value = "★bot ◆ツール"

# Row 4
$ws.Range("A4").Value = "2025-12-11 06:30:42"
$ws.Range("B4").Value = "自動出品システムの開発"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5451514"
$ws.Range("G4").Value = 83
$ws.Range("H4").Value = "◆開発"

# Row 5
$ws.Range("A5").Value = "2025-12-11 06:30:42"
$ws.Range("B5").Value = "Access DB家賃管理SYSを最新Access で稼働できるように"
$ws.Range("C5").Value = "システム開発"
$ws.Range("D5").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E5").Value = "期限情報なし"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5451626"
$ws.Range("G5").Value = 38
$ws.Range("H5").Value = "◇管理"

# Row 6
$ws.Range("A6").Value = "2025-12-11 06:30:42"
$ws.Range("B6").Value = "進行管理およびチームディレクションを担当"
$ws.Range("C6").Value = "システム開発"
$ws.Range("D6").Value = "~ 5,000 円 / 固定"
$ws.Range("E6").Value = "期限情報なし"
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5418064"
$ws.Range("G6").Value = 30
$ws.Range("H6").Value = "◇管理"

# Row 7
$ws.Range("A7").Value = "2025-12-11 06:30:42"
$ws.Range("B7").Value = "Rubyの暗号化機能のPHP化"
$ws.Range("C7").Value = "システム開発"
$ws.Range("D7").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E7").Value = "期限情報なし"
$ws.Range("F7").Value = "https://www.lancers.jp/work/detail/5451714"
$ws.Range("G7").Value = 28
$ws.Range("H7").Value = "○PHP"

# Row 8
$ws.Range("A8").Value = "2025-12-11 06:30:42"
$ws.Range("B8").Value = "【オンライン講師募集】バックエンドの基礎を教えていただける方"
$ws.Range("C8").Value = "システム開発"
$ws.Range("D8").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E8").Value = "期限情報なし"
$ws.Range("F8").Value = "https://www.lancers.jp/work/detail/5451420"
$ws.Range("G8").Value = 18
$ws.Range("H8").ClearContents()

# Row 9
$ws.Range("A9").Value = "2025-12-11 06:30:42"
$ws.Range("B9").Value = "注目 限定公開 PR 限定公開の仕事"
$ws.Range("C9").Value = "システム開発"
$ws.Range("D9").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E9").Value = "期限情報なし"
$ws.Range("F9").Value = "https://www.lancers.jp/work/detail/5450323"
$ws.Range("G9").Value = 13
$ws.Range("H9").ClearContents()
